# Update the cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for the rows whose figures changed, and swap the TRON / BinanceUSD
# rows (16 and 17) so TRON now appears before BinanceUSD with refreshed data.
#
# Column D holds price figures formatted as plain text (e.g. "29.416.00",
# "1.009", "91.14"); force text format before writing so Excel doesn't
# reinterpret values such as "1.007" as a numeric literal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    # Column D prices are plain text (e.g. "29.416.00", "1.009", "91.14").
    # Values that use a European thousands-separator style (more than one
    # '.') can never round-trip through Excel's numeric parser, so they are
    # safe to assign directly. Values that parse as a plain number (e.g.
    # "1.007") would otherwise be silently reinterpreted as a numeric
    # literal, so force text format first for just those cells.
    param($Row, $Text)
    $cell = $ws.Cells.Item($Row, 4)
    if ($Text -match '^[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $Text
}

function Set-VolumeText {
    param($Row, $Text)
    $ws.Cells.Item($Row, 5).Value = $Text
}

# Row 16 becomes TRON, row 17 becomes BinanceUSD (swap of the two rows'
# Coin name, Link, Price and Volume(1h) values).
$ws.Cells.Item(16, 2).Value = "TRON"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-PriceText 16 "0.06802"
Set-VolumeText 16 "  +2.45%  "

$ws.Cells.Item(17, 2).Value = "BinanceUSD"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-PriceText 17 "1.009"
Set-VolumeText 17 "  +0.84%  "

# Price/Volume refreshes for the remaining rows.
$updates = @(
    @{ Row = 2;  Price = "29.391.22";   Volume = "  +0.57%  " },
    @{ Row = 3;  Price = "1.911.99";    Volume = "  +1.06%  " },
    @{ Row = 4;  Price = "1.007";       Volume = "  +0.60%  " },
    @{ Row = 5;  Price = "325.15";      Volume = "  +0.65%  " },
    @{ Row = 6;  Price = $null;         Volume = "  +0.73%  " },
    @{ Row = 7;  Price = "0.4818";      Volume = "  +1.37%  " },
    @{ Row = 8;  Price = "0.4065";      Volume = "  +0.51%  " },
    @{ Row = 9;  Price = "0.08223";     Volume = "  +2.43%  " },
    @{ Row = 10; Price = $null;         Volume = "  +2.27%  " },
    @{ Row = 11; Price = $null;         Volume = "  +1.40%  " },
    @{ Row = 12; Price = "1.934.83";    Volume = "  +2.88%  " },
    @{ Row = 13; Price = "6.035";       Volume = "  +1.98%  " },
    @{ Row = 14; Price = "7.207";       Volume = "  +2.58%  " },
    @{ Row = 15; Price = "91.09";       Volume = "  +2.23%  " },
    @{ Row = 18; Price = "0.00001038";  Volume = "  +0.90%  " },
    @{ Row = 19; Price = "17.70";       Volume = "  +1.03%  " },
    @{ Row = 20; Price = $null;         Volume = "  +0.66%  " },
    @{ Row = 21; Price = "29.438.10";   Volume = "  +0.72%  " },
    @{ Row = 22; Price = "5.626";       Volume = "  +2.04%  " },
    @{ Row = 23; Price = "11.80";       Volume = "  +1.25%  " },
    @{ Row = 24; Price = "2.190";       Volume = "  +1.42%  " },
    @{ Row = 25; Price = "2.117.35";    Volume = "  +0.57%  " },
    @{ Row = 26; Price = "6.557";       Volume = "  +10.73%  " },
    @{ Row = 27; Price = "156.57";      Volume = "  +1.51%  " },
    @{ Row = 28; Price = "20.00";       Volume = "  +1.45%  " },
    @{ Row = 29; Price = "2.103";       Volume = "  +0.96%  " },
    @{ Row = 30; Price = "120.10";      Volume = "  +2.01%  " },
    @{ Row = 31; Price = "1.017";       Volume = "  -0.56%  " },
    @{ Row = 32; Price = "0.09548";     Volume = "  +1.32%  " },
    @{ Row = 33; Price = "5.582";       Volume = "  +4.61%  " },
    @{ Row = 34; Price = "3.550";       Volume = "  +0.63%  " },
    @{ Row = 35; Price = "1.364";       Volume = "  -0.57%  " },
    @{ Row = 36; Price = "0.02282";     Volume = "  +1.63%  " },
    @{ Row = 37; Price = "0.06113";     Volume = "  +1.39%  " },
    @{ Row = 38; Price = "1.179";       Volume = "  +1.06%  " },
    @{ Row = 39; Price = "8.047";       Volume = "  +2.39%  " },
    @{ Row = 40; Price = "0.5961";      Volume = "  +2.33%  " },
    @{ Row = 41; Price = "10.81";       Volume = "  +7.97%  " },
    @{ Row = 42; Price = "0.1847";      Volume = "  +0.93%  " },
    @{ Row = 43; Price = $null;         Volume = "  -0.50%  " },
    @{ Row = 44; Price = "2.384";       Volume = "  +1.39%  " },
    @{ Row = 45; Price = "0.07611";     Volume = "  -1.10%  " },
    @{ Row = 46; Price = "12.37";       Volume = "  +2.28%  " },
    @{ Row = 47; Price = "0.5567";      Volume = "  +1.57%  " },
    @{ Row = 48; Price = "1.951";       Volume = "  +2.32%  " },
    @{ Row = 49; Price = "117.65";      Volume = "  +4.28%  " },
    @{ Row = 50; Price = "2.427";       Volume = "  +4.40%  " },
    @{ Row = 51; Price = $null;         Volume = "  +1.05%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.Price) {
        Set-PriceText $u.Row $u.Price
    }
    Set-VolumeText $u.Row $u.Volume
}
